# Applies the "blenda script" textbox tweak (position nudge + run split)
# to the "Google Shape;85;p1" textbox on every slide of the template deck.
#
# Position is stored in the XML as EMU (1 pt = 12700 EMU), but the
# PowerPoint Shape object exposes Left/Top/Width/Height in points as a
# Single (32-bit float), so the literals below are chosen so that they
# round-trip through Single precision to the exact target EMU values:
#   Left   -> 5017590 EMU
#   Top    -> 4699823 EMU
#   Width  -> 2794500 EMU (unchanged)
#   Height -> 400200  EMU (unchanged)

$targetLeft   = 395.08586614173225
$targetTop    = 370.06484251968504
$targetWidth  = 220.0394094488189
$targetHeight = 31.51185039370079

$p = $ppt.ActivePresentation

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $s = $p.Slides.Item($slideIdx)

    $shape = $null
    foreach ($candidate in $s.Shapes) {
        if ($candidate.Name -eq "Google Shape;85;p1") {
            $shape = $candidate
            break
        }
    }
    if ($shape -eq $null) { continue }

    $tr = $shape.TextFrame.TextRange
    if ($tr.Text -ne "blenda script") { continue }

    # Split the single "blenda script" run into "blenda" + " script" runs
    # (matching the source diff) by touching each sub-range's font.
    $first = $tr.Characters(1, 6)
    $first.Font.Name = "Calibri"

    $second = $tr.Characters(7, 7)
    $second.Font.Name = "Calibri"

    # Re-assert the on-screen position/size: editing the text can trigger
    # an auto-fit re-layout, so restore the exact EMU geometry afterwards.
    $shape.Left = $targetLeft
    $shape.Top = $targetTop
    $shape.Width = $targetWidth
    $shape.Height = $targetHeight
}
